# Fix bugs in VAR and ARIMA results across the four sheets.
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet 1: y_fitted_on_begin_2016
# Insert a new row at row 2 (year 1991) and rewrite all y_value data.
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws1.Rows.Item(2).Insert()

$sheet1Data = @(
    @(1991, 61.43048008274011),
    @(1992, 61.78157261392445),
    @(1993, 62.19834379772002),
    @(1994, 62.63976656818687),
    @(1995, 62.92423120383471),
    @(1996, 63.06875156727907),
    @(1997, 63.25185562200105),
    @(1998, 63.54078961863272),
    @(1999, 63.72286240574709),
    @(2000, 63.92143206163952),
    @(2001, 64.09951070100728),
    @(2002, 63.99906582262692),
    @(2003, 64.08261901273161),
    @(2004, 64.60517168785685),
    @(2005, 64.55096774851449),
    @(2006, 64.18026736846923),
    @(2007, 64.29180032272009),
    @(2008, 64.53953370835438),
    @(2009, 64.3277694120067),
    @(2010, 64.59289732081527),
    @(2011, 64.69323366409317),
    @(2012, 65.29724755797528),
    @(2013, 66.19284351606511),
    @(2014, 66.63967955069324),
    @(2015, 67.00636838519286),
    @(2016, 67.87289279891756)
)

$r = 2
foreach ($row in $sheet1Data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# -----------------------------------------------------------------
# Sheet 2: y_pred_on_2017_2021
# Same number of rows (2017-2021); only y_value column changes.
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("y_pred_on_2017_2021")

$sheet2Values = @(69.66185396507895, 71.33462352133481, 73.25006698724667, 75.71730684921548, 78.80762897570685)

$r = 2
foreach ($v in $sheet2Values) {
    $ws2.Cells.Item($r, 2).Value = $v
    $r = $r + 1
}

# -----------------------------------------------------------------
# Sheet 3: y_fitted_on_begin_2021
# Delete row 2 (year 1990) and rewrite all remaining y_value data.
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws3.Rows.Item(2).Delete()

$sheet3Data = @(
    @(1991, 61.0856468906162),
    @(1992, 61.89954979427429),
    @(1993, 62.42014308655836),
    @(1994, 62.9022968424083),
    @(1995, 63.00133272773329),
    @(1996, 63.22672064413919),
    @(1997, 62.9198800439135),
    @(1998, 63.49757689768331),
    @(1999, 63.83299998445116),
    @(2000, 64.29730162145609),
    @(2001, 64.38847888834567),
    @(2002, 64.04653692164901),
    @(2003, 64.18347929467224),
    @(2004, 64.2486166982343),
    @(2005, 64.95937864210968),
    @(2006, 64.21436323880738),
    @(2007, 64.29496908652393),
    @(2008, 65.01712382616167),
    @(2009, 64.9637891040161),
    @(2010, 64.43784233847023),
    @(2011, 65.2118970856898),
    @(2012, 64.4123616933975),
    @(2013, 65.66047326822748),
    @(2014, 66.06958099384612),
    @(2015, 65.85075027324396),
    @(2016, 67.28037443138376),
    @(2017, 67.59463799184418),
    @(2018, 66.8698106703667),
    @(2019, 66.77268622663917),
    @(2020, 65.76554142229169),
    @(2021, 65.85691964904467)
)

$r = 2
foreach ($row in $sheet3Data) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# -----------------------------------------------------------------
# Sheet 4: y_pred_on_2022_2026
# Same number of rows (2022-2026); only y_value column changes.
# -----------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("y_pred_on_2022_2026")

$sheet4Values = @(64.61813990753957, 63.88571096701131, 63.24741749443844, 62.60519417815841, 61.94577354948412)

$r = 2
foreach ($v in $sheet4Values) {
    $ws4.Cells.Item($r, 2).Value = $v
    $r = $r + 1
}
